$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.384.62'
$ws.Range("E2").Value = '  +4.22%  '

$ws.Range("D3").Value = '2.428.25'
$ws.Range("E3").Value = '  +5.58%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = "'556.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.92%  '

$ws.Range("D6").Value = "'138.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.21%  '

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").Value = "'0.585"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.02%  '

$ws.Range("D9").Value = '2.426.58'
$ws.Range("E9").Value = '  +5.50%  '

$ws.Range("D10").Value = "'0.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.80%  '

$ws.Range("D11").Value = "'5.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.32%  '

$ws.Range("D12").Value = "'0.151"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.49%  '

$ws.Range("D13").Value = "'0.348"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.18%  '

$ws.Range("D14").Value = "'26.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +13.10%  '

$ws.Range("D15").Value = '2.863.73'
$ws.Range("E15").Value = '  +5.81%  '

$ws.Range("D16").Value = '62.306.19'
$ws.Range("E16").Value = '  +4.18%  '

$ws.Range("D17").Value = "'0.0000141"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.00%  '

$ws.Range("D18").Value = '2.434.39'
$ws.Range("E18").Value = '  +5.70%  '

$ws.Range("D19").Value = "'11.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.16%  '

$ws.Range("D20").Value = "'346.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.17%  '

$ws.Range("D21").Value = "'4.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.17%  '

$ws.Range("D22").Value = "'6.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.00%  '

$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("D24").Value = "'65.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.31%  '

$ws.Range("E25").Value = '  +1.63%  '

$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("D27").Value = "'1.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.79%  '

$ws.Range("D28").Value = "'8.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.83%  '

$ws.Range("E29").Value = '  +14.53%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'1.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.67%  '

$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '0.0₃0783'
$ws.Range("E31").Value = '  +8.40%  '

$ws.Range("D32").Value = "'6.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.39%  '

$ws.Range("D33").Value = "'171.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.79%  '

$ws.Range("E34").Value = '  +6.61%  '

$ws.Range("E35").Value = '  +5.50%  '

$ws.Range("D36").Value = "'380.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +20.11%  '

$ws.Range("D37").Value = "'18.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.24%  '

$ws.Range("D38").Value = "'4.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.93%  '

$ws.Range("E39").Value = '  -0.03%  '

$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.09%  '

$ws.Range("E41").Value = '  +12.27%  '

$ws.Range("D42").Value = "'39.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.56%  '

$ws.Range("D43").Value = "'145.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.70%  '

$ws.Range("D44").Value = "'3.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.07%  '

$ws.Range("D45").Value = "'20.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +10.61%  '

$ws.Range("D46").Value = "'0.590"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.98%  '

$ws.Range("D47").Value = "'0.0953"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.92%  '

$ws.Range("D48").Value = "'0.0519"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.21%  '

$ws.Range("D49").Value = "'0.0222"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.28%  '

$ws.Range("D50").Value = "'17.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.05%  '

$ws.Range("D51").Value = '0.0₆0218'
$ws.Range("E51").Value = '  -2.03%  '
